$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Clean Code A Handbook of Agile Software Craftsmanship"
$ws.Range("C2").Value = "['James 0. Coplien', 'Robert C Martin']"
$ws.Range("D2").Value = "['PReNtICE HALL']"
